$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing header style (cellXfs index 1, bold/centered/bordered)
# by copying A1's formatting to the new header cells (B1:F1) and to the new
# A-column id cells (A2:A5) BEFORE we overwrite any values.
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 1 (header): now B=lang_code, C=code, D=name, E=descr, F=is_active.
# A1 no longer holds a header label; drop the cell entirely.
$ws.Range("A1").Clear()
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eng"
$ws.Range("C2").Value = "MNA"
$ws.Range("D2").Value = "Manual Adjudication"
$ws.Range("E2").Value = "Rejection during Manual Adjudication"
$ws.Range("F2").Value = $true

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eng"
$ws.Range("C3").Value = "CLR"
$ws.Range("D3").Value = "Client Rejection"
$ws.Range("E3").Value = "Rejection in Registration Client"
$ws.Range("F3").Value = $true

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "fra"
$ws.Range("C4").Value = "MNA"
$ws.Range("D4").Value = "Manuel arbitrage"
$ws.Range("E4").Value = "Renvoi en cours de sélection manuelle"
$ws.Range("F4").Value = $true

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "fra"
$ws.Range("C5").Value = "CLR"
$ws.Range("D5").Value = "Rejet de client"
$ws.Range("E5").Value = "Rejet en enregistrement Client"
$ws.Range("F5").Value = $true
